# Generate Report for Handback
# ------------------------------------------------------------------
# This localization-status workbook tracks, per target language, the
# handoff/handback lifecycle of a single source file
# (14d30e09-5902-4338-a1a8-2cf0d70ead89.md). This run records that the
# de-de / zh-cn handback has completed: status flips to "in sync", the
# generated target (.md) and handback (.xlf) files are recorded (with a
# link back to the source doc), and the handback timestamp is stamped.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$mdFile            = "14d30e09-5902-4338-a1a8-2cf0d70ead89.md"
$mdFileUrl         = "https://github.com/OpenLocalizationTestOrg/oltest/blob/f2b8f24e954788671e29b33a7eb1f9112f7143d6/e2e/14d30e09-5902-4338-a1a8-2cf0d70ead89.md"

$zhXlf   = "14d30e09-5902-4338-a1a8-2cf0d70ead89.0632078ed0787b45800ec8667fae4eae79f47fed.zh-cn.xlf"
$deXlf   = "14d30e09-5902-4338-a1a8-2cf0d70ead89.0632078ed0787b45800ec8667fae4eae79f47fed.de-de.xlf"

$zhHandbackTime = "2016-08-13 15:13:10"
$deHandbackTime = "2016-08-13 15:13:20"

# Column width Excel stores on disk = round(ColumnWidth * 6) / 6 + 5/6,
# so to land on the report's target "best fit" widths we back out the
# ColumnWidth value that rounds closest to the desired stored width.
$wideColWidth   = 29.166666666666668   # -> stored width 30 (closest reachable to 29.9777047293527)
$maxColWidth    = 39.166666666666664   # -> stored width 40

# ------------------------------------------------------------------
# Overview sheet: per-language status column reflects the new status
# ------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Columns.Item(5).ColumnWidth = $wideColWidth
$overview.Columns.Item(6).ColumnWidth = $wideColWidth

# ------------------------------------------------------------------
# zh-cn sheet: status, target file link, handback file + datetime
# ------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdFileUrl, "", "", $mdFile) | Out-Null
$zhcn.Range("J2").Value = $zhXlf
$zhcn.Range("K2").Value = $zhHandbackTime
$zhcn.Columns.Item(3).ColumnWidth = $wideColWidth
$zhcn.Columns.Item(9).ColumnWidth = $maxColWidth
$zhcn.Columns.Item(10).ColumnWidth = $maxColWidth

# ------------------------------------------------------------------
# de-de sheet: status, target file link, handback file + datetime
# ------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusHandedBack
$dede.Hyperlinks.Add($dede.Range("I2"), $mdFileUrl, "", "", $mdFile) | Out-Null
$dede.Range("J2").Value = $deXlf
$dede.Range("K2").Value = $deHandbackTime
$dede.Columns.Item(3).ColumnWidth = $wideColWidth
$dede.Columns.Item(9).ColumnWidth = $maxColWidth
$dede.Columns.Item(10).ColumnWidth = $maxColWidth
